$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Descriptions (column D) for the new rows first
$ws.Range("D5").Value = "The fourth test parameter"
$ws.Range("D6").Value = "The fifth test parameter"
$ws.Range("D7").Value = "The sixth test parameter"

# Names (column C) for the new rows
$ws.Range("C5").Value = "Test Parameter 4"
$ws.Range("C6").Value = "Test Parameter 5"
$ws.Range("C7").Value = "Test Parameter 6"

# Units (column E) for the new rows
$ws.Range("E5").Value = "Some Units"
$ws.Range("E6").Value = "Some Units"
$ws.Range("E7").Value = "Some Units"

# Distribution types (column F) for the new rows
$ws.Range("F5").Value = "Log Uniform"
$ws.Range("F6").Value = "Truncated Normal"
$ws.Range("F7").Value = "Truncated Log Normal"

# Match the author's final active selection after entering the new data
$ws.Range("D14").Select() | Out-Null
